$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (rows 2-18), columns A=Player, B=Position, C=Team
$data = @(
    @("Stephen Curry",      "PG",    "Golden State Warriors"),
    @("Malcolm Brogdon",    "PG,SG", "Washington Wizards"),
    @("Austin Reaves",      "PG,SG", "Los Angeles Lakers"),
    @("Darius Garland",     "PG",    "Cleveland Cavaliers"),
    @("OG Anunoby",         "SF,PF", "New York Knicks"),
    @("Kevin Durant",       "SF,PF", "Phoenix Suns"),
    @("Karl-Anthony Towns", "PF,C",  "New York Knicks"),
    @("Trey Murphy III",    "SF,PF", "New Orleans Pelicans"),
    @("Mark Williams",      "C",     "Charlotte Hornets"),
    @("Jarrett Allen",      "C",     "Cleveland Cavaliers"),
    @("Jalen Duren",        "C",     "Detroit Pistons"),
    @("Keegan Murray",      "SF,PF", "Sacramento Kings"),
    @("Tyrese Maxey",       "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton",  "PG,SG", "Indiana Pacers"),
    @("Daniel Gafford",     "PF,C",  "Dallas Mavericks"),
    @("Franz Wagner",       "SF,PF", "Orlando Magic"),
    @("Jalen Johnson",      "SF,PF", "Atlanta Hawks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
